# Update "Elenco attività per i progetti" - aggiornamento tabella di marcia
# Updates the "% completamento" values for a few phases on the "Progetto 1" sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Progetto 1")

# Row 12 - "Pianificazione Web Server": 0.05 -> 0.15
$ws.Range("B12").Value = 0.15

# Row 13 - "Gestione delle query del server": 0 -> 0.35
$ws.Range("B13").Value = 0.35

# Row 15 - "Gestione connessione al DataBase": 0 -> 1
$ws.Range("B15").Value = 1

# Update the selected cell on the sheet to match the saved view state
$ws.Activate()
$ws.Range("B11").Select()
